$d = $word.ActiveDocument

# --- Step 1: Paragraph 1 ("Testgit ! ") -----------------------------------
# The original paragraph holds two runs ("Testgit" and " ! ") wrapped by a
# <w:proofErr> spell-check pair. Deleting the whole paragraph (mark
# included) drops the stray proofErr bookkeeping, then we re-create a
# clean paragraph ahead of "Hello world." with a single merged run.
$pFirst = $d.Paragraphs(1)
$pFirst.Range.Delete()

$pHello = $d.Paragraphs(1)
$pHello.Range.InsertParagraphBefore()
$pNew = $d.Paragraphs(1)
$pNew.Range.Text = "Testgit ! "

# --- Step 2: add the new third paragraph after "Hello world." -------------
$pHello = $d.Paragraphs(2)
$pHello.Range.InsertParagraphAfter()
$pThird = $d.Paragraphs(3)

$finalText = "C" + [char]0x2019 + "est une tentative de modification de fichier."
# Append a placeholder trailing character first: this keeps the bookmark
# insertion point (added below) away from the paragraph-mark slot, which
# this COM host mis-resolves to the very start of the document when a
# collapsed Range sits exactly on a paragraph mark.
$pThird.Range.Text = $finalText + "X"

# --- Step 3: move the _GoBack bookmark onto the new last paragraph --------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$bmPos = $pThird.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the placeholder character now that the bookmark is anchored.
$placeholder = $d.Range($pThird.Range.End - 2, $pThird.Range.End - 1)
$placeholder.Delete()

Write-Output "ok"
